$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = -2

# Update the selected/active cell in the bottom-right frozen pane
$ws.Range("D7").Select()
